$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("D1").Value = "Shape (obs x variables)"
$ws.Range("G1").Value = "Runtime "

# --- Row 3 updates (existing MICE default-method run becomes the combined MICE + Demographics row) ---
$ws.Range("B3").Value = "MICE Imputation  on NHANES Questionnaire Data  + Demographics"

# --- Row 4 updates (existing "cart" run row gets corrected numbers / new note, plus new columns) ---
$ws.Range("B4").Value = 'MICE Imputation  on NHANES Questionnaire Data (Diet Behavior & Quality + Demographics), method = "cart" (classification & Linear regression)'
$ws.Range("D4").Value = "9971 x 97"
$ws.Range("E4").Value = 'data <- mice(data, m = 1, seed = 2022, method = "cart")'
$ws.Range("F4").Value = "Yes"
$ws.Range("G4").Value = "6.328955 min"
$ws.Range("H4").Value = "1. default method is stochastic, while 'cart' is not 2. # of Logged events: 176, 3. WRONG DATASET USED (discovered after running)"
$ws.Range("I4").Value = "None"

# row 4 grows taller to fit the expanded note text
$ws.Rows(4).RowHeight = 72

# --- New row 5: re-run of the "cart" MICE imputation on the corrected dataset ---
$ws.Range("A5").Value = 44575
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats -> reuse the date number format (style 1)

$ws.Range("B5").Value = 'MICE Imputation  on NHANES Questionnaire Data (Diet Behavior & Quality + Demographics), method = "cart" (classification & Linear regression)'
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)   # reuse wrap-text style (style 2)

$ws.Range("C5").Value = "2017-18 NHANES Questionnaire Data "
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)   # reuse wrap-text style (style 2)

$ws.Range("D5").Value = "9254 x 91"
$ws.Range("E5").Value = 'data <- mice(data, m = 1, seed = 2022, method = "cart)'
$ws.Range("F5").Value = "Yes"
$ws.Range("G5").Value = "22.45882 min"

$ws.Range("H5").Value = "# of Logged Events: 587"
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)   # reuse wrap-text style (style 2)

$ws.Rows(5).RowHeight = 43.2

# --- New row 6: placeholder row for a further run (only Date + Test populated so far) ---
$ws.Range("A6").Value = 44575
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # reuse the date number format (style 1)

$ws.Range("B6").Value = 'MICE Imputation  on NHANES Questionnaire Data (Diet Behavior & Quality + Demographics), method = "cart" (classification & Linear regression)'
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)   # reuse wrap-text style (style 2)

$ws.Rows(6).RowHeight = 43.2

$excel.CutCopyMode = $false

# --- View state: selection moves to D1, window scrolled so column C is leftmost ---
$excel.Goto($ws.Range("C1"), $true)
$ws.Range("D1").Select()
